$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.189005374908447
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 2343.969756180612
$ws.Range("F2").Value = 0.1379218287815959
$ws.Range("G2").Value = 0.09762381338770172
$ws.Range("H2").Value = 0.08392823345613033
$ws.Range("I2").Value = 0.07541723828400335
$ws.Range("J2").Value = 0.06718826145539383
$ws.Range("K2").Value = 0.06494415141123731
$ws.Range("L2").Value = 0.06033705017073068
$ws.Range("M2").Value = 0.05526993044324594
$ws.Range("N2").Value = 0.05382811397613153
$ws.Range("O2").Value = 0.05159656540095481
$ws.Range("P2").Value = 0.05019217888448748
$ws.Range("Q2").Value = 0.04903716297000969
$ws.Range("R2").Value = 0.04771894309870181
$ws.Range("S2").Value = 0.04720876430970275
$ws.Range("T2").Value = 0.04720876430970275
$ws.Range("U2").Value = 0.0469972698320686
$ws.Range("V2").Value = 0.04604861009825346
$ws.Range("W2").Value = 0.04604861009825346
$ws.Range("X2").Value = 0.04586526025393725
$ws.Range("Y2").Value = 0.04569141824913472

$ws.Range("C3").Value = 1.176010608673096
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 2369.482337174024
$ws.Range("F3").Value = 0.1487985439548758
$ws.Range("G3").Value = 0.1078401265590551
$ws.Range("H3").Value = 0.0896557191526655
$ws.Range("I3").Value = 0.07501903566353273
$ws.Range("J3").Value = 0.06546581566684902
$ws.Range("K3").Value = 0.05974707609369931
$ws.Range("L3").Value = 0.05738712223894216
$ws.Range("M3").Value = 0.05451969127793981
$ws.Range("N3").Value = 0.05174570081600877
$ws.Range("O3").Value = 0.05107067444195169
$ws.Range("P3").Value = 0.0505421992273908
$ws.Range("Q3").Value = 0.04978642770906904
$ws.Range("R3").Value = 0.04876478165922864
$ws.Range("S3").Value = 0.0481396497895012
$ws.Range("T3").Value = 0.04784962740959437
$ws.Range("U3").Value = 0.04727141415030775
$ws.Range("V3").Value = 0.04666400869139194
$ws.Range("W3").Value = 0.0466555972452058
$ws.Range("X3").Value = 0.04640179076413146
$ws.Range("Y3").Value = 0.04618873951606284

$ws.Range("C4").Value = 1.269996404647827
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 2325.246779653955
$ws.Range("F4").Value = 0.1354892594373364
$ws.Range("G4").Value = 0.1043364267237635
$ws.Range("H4").Value = 0.08641855357299737
$ws.Range("I4").Value = 0.07233382944655538
$ws.Range("J4").Value = 0.06409622841075759
$ws.Range("K4").Value = 0.06121366107722484
$ws.Range("L4").Value = 0.05524528878660095
$ws.Range("M4").Value = 0.05270334929793571
$ws.Range("N4").Value = 0.05077378353248418
$ws.Range("O4").Value = 0.04918864993969208
$ws.Range("P4").Value = 0.04854822603863175
$ws.Range("Q4").Value = 0.04822649714705302
$ws.Range("R4").Value = 0.04751396050366662
$ws.Range("S4").Value = 0.04686705761487967
$ws.Range("T4").Value = 0.04639694252760584
$ws.Range("U4").Value = 0.04612762504233548
$ws.Range("V4").Value = 0.04575788202699366
$ws.Range("W4").Value = 0.04554014773318141
$ws.Range("X4").Value = 0.0453596445210215
$ws.Range("Y4").Value = 0.04532644794647085

$ws.Range("C5").Value = 1.308999061584473
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 2525.808780989515
$ws.Range("F5").Value = 0.1480006208173046
$ws.Range("G5").Value = 0.105940132348289
$ws.Range("H5").Value = 0.08692557005632405
$ws.Range("I5").Value = 0.08148431452949224
$ws.Range("J5").Value = 0.0703028626641337
$ws.Range("K5").Value = 0.06628566010386341
$ws.Range("L5").Value = 0.06117644800618051
$ws.Range("M5").Value = 0.05880804561412753
$ws.Range("N5").Value = 0.05542723465798415
$ws.Range("O5").Value = 0.05409100981600161
$ws.Range("P5").Value = 0.05306805376483209
$ws.Range("Q5").Value = 0.05230005403519338
$ws.Range("R5").Value = 0.05140924370025027
$ws.Range("S5").Value = 0.05085342462663037
$ws.Range("T5").Value = 0.05071901371862765
$ws.Range("U5").Value = 0.05017348592765
$ws.Range("V5").Value = 0.04978694260503969
$ws.Range("W5").Value = 0.04962833721928785
$ws.Range("X5").Value = 0.0494483405729781
$ws.Range("Y5").Value = 0.04923603861578001

$ws.Range("C6").Value = 1.155039310455322
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 2316.943882444952
$ws.Range("F6").Value = 0.1403179557676519
$ws.Range("G6").Value = 0.09620524580068289
$ws.Range("H6").Value = 0.08126550927451952
$ws.Range("I6").Value = 0.07122758190643134
$ws.Range("J6").Value = 0.06311055976925004
$ws.Range("K6").Value = 0.05753835023355077
$ws.Range("L6").Value = 0.05297124610263465
$ws.Range("M6").Value = 0.05187066600338757
$ws.Range("N6").Value = 0.04928243544552424
$ws.Range("O6").Value = 0.04898700179135378
$ws.Range("P6").Value = 0.04846039758784792
$ws.Range("Q6").Value = 0.04759546871124062
$ws.Range("R6").Value = 0.04717171223692289
$ws.Range("S6").Value = 0.04677193801497184
$ws.Range("T6").Value = 0.04636800158069497
$ws.Range("U6").Value = 0.04597416131154718
$ws.Range("V6").Value = 0.04558225669692159
$ws.Range("W6").Value = 0.04547224971162065
$ws.Range("X6").Value = 0.04521394302413235
$ws.Range("Y6").Value = 0.04516459809834213

$ws.Range("C7").Value = 1.179968118667603
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 2345.202352554741
$ws.Range("F7").Value = 0.1417727471006475
$ws.Range("G7").Value = 0.1017164400338181
$ws.Range("H7").Value = 0.08315092676843923
$ws.Range("I7").Value = 0.07770210117651227
$ws.Range("J7").Value = 0.07108003110657272
$ws.Range("K7").Value = 0.06459266102253855
$ws.Range("L7").Value = 0.05641355419285166
$ws.Range("M7").Value = 0.05587372371434392
$ws.Range("N7").Value = 0.05337208073408346
$ws.Range("O7").Value = 0.05167963185511477
$ws.Range("P7").Value = 0.04914248349112386
$ws.Range("Q7").Value = 0.04761306877921351
$ws.Range("R7").Value = 0.04714499486170638
$ws.Range("S7").Value = 0.04682895599613247
$ws.Range("T7").Value = 0.04643530890870257
$ws.Range("U7").Value = 0.04633943309342478
$ws.Range("V7").Value = 0.04585158143246637
$ws.Range("W7").Value = 0.04578482721956155
$ws.Range("X7").Value = 0.04571544546890333
$ws.Range("Y7").Value = 0.04571544546890333

$ws.Range("C8").Value = 1.250036239624023
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 2333.872984075672
$ws.Range("F8").Value = 0.1490791621638265
$ws.Range("G8").Value = 0.1011527852661363
$ws.Range("H8").Value = 0.08720701755267757
$ws.Range("I8").Value = 0.07474866858658949
$ws.Range("J8").Value = 0.06735947790126152
$ws.Range("K8").Value = 0.06149695676410665
$ws.Range("L8").Value = 0.05815731382296782
$ws.Range("M8").Value = 0.05631878535785891
$ws.Range("N8").Value = 0.05370350022005747
$ws.Range("O8").Value = 0.05182159117701591
$ws.Range("P8").Value = 0.04902860776769753
$ws.Range("Q8").Value = 0.04867343726322359
$ws.Range("R8").Value = 0.04809926703088892
$ws.Range("S8").Value = 0.04745364331230157
$ws.Range("T8").Value = 0.04728651629212414
$ws.Range("U8").Value = 0.04698256215715222
$ws.Range("V8").Value = 0.04648572228080717
$ws.Range("W8").Value = 0.04597237013578975
$ws.Range("X8").Value = 0.04564339219947164
$ws.Range("Y8").Value = 0.04549460007944779

$ws.Range("C9").Value = 1.064001560211182
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 2366.507276012098
$ws.Range("F9").Value = 0.1454388795598044
$ws.Range("G9").Value = 0.1051601298381576
$ws.Range("H9").Value = 0.08597291051193648
$ws.Range("I9").Value = 0.07409327504012035
$ws.Range("J9").Value = 0.06650377873874121
$ws.Range("K9").Value = 0.06155401438176421
$ws.Range("L9").Value = 0.05554200399772167
$ws.Range("M9").Value = 0.05243088377181705
$ws.Range("N9").Value = 0.0522577203298514
$ws.Range("O9").Value = 0.05018152935252544
$ws.Range("P9").Value = 0.04931724771792385
$ws.Range("Q9").Value = 0.04816939541101047
$ws.Range("R9").Value = 0.04816939541101047
$ws.Range("S9").Value = 0.04753367241670278
$ws.Range("T9").Value = 0.0467229329537231
$ws.Range("U9").Value = 0.0467229329537231
$ws.Range("V9").Value = 0.04654715543157777
$ws.Range("W9").Value = 0.04627330289001558
$ws.Range("X9").Value = 0.04627330289001558
$ws.Range("Y9").Value = 0.04613074612109352

$ws.Range("C10").Value = 1.258014678955078
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 2375.033445383437
$ws.Range("F10").Value = 0.1244215905629234
$ws.Range("G10").Value = 0.1015054022390827
$ws.Range("H10").Value = 0.08328388809532412
$ws.Range("I10").Value = 0.07283539889984959
$ws.Range("J10").Value = 0.06843615063589194
$ws.Range("K10").Value = 0.06186319894778834
$ws.Range("L10").Value = 0.05764173825951959
$ws.Range("M10").Value = 0.05583841212785171
$ws.Range("N10").Value = 0.05453386465467414
$ws.Range("O10").Value = 0.05209776043278985
$ws.Range("P10").Value = 0.05048826401231109
$ws.Range("Q10").Value = 0.04966632133697318
$ws.Range("R10").Value = 0.04871323346025792
$ws.Range("S10").Value = 0.0478448239162141
$ws.Range("T10").Value = 0.04745204438538783
$ws.Range("U10").Value = 0.04745204438538783
$ws.Range("V10").Value = 0.04701371624263587
$ws.Range("W10").Value = 0.04655213535004091
$ws.Range("X10").Value = 0.04637328657262785
$ws.Range("Y10").Value = 0.04629694825308842

$ws.Range("C11").Value = 1.124999046325684
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 2396.645271401234
$ws.Range("F11").Value = 0.1356765322885027
$ws.Range("G11").Value = 0.1051430614231803
$ws.Range("H11").Value = 0.09024627092258045
$ws.Range("I11").Value = 0.08039669197952146
$ws.Range("J11").Value = 0.07372249716436251
$ws.Range("K11").Value = 0.06686882904582876
$ws.Range("L11").Value = 0.06318789564759339
$ws.Range("M11").Value = 0.06057242971060489
$ws.Range("N11").Value = 0.05858379284136753
$ws.Range("O11").Value = 0.05641980513892528
$ws.Range("P11").Value = 0.05366881585895579
$ws.Range("Q11").Value = 0.05150677702919379
$ws.Range("R11").Value = 0.05148920369848264
$ws.Range("S11").Value = 0.0504528867205261
$ws.Range("T11").Value = 0.04932340616522807
$ws.Range("U11").Value = 0.04864520748761079
$ws.Range("V11").Value = 0.04743061002416035
$ws.Range("W11").Value = 0.04706947534926163
$ws.Range("X11").Value = 0.04682790697688306
$ws.Range("Y11").Value = 0.04671823141133009

